# [ADD] initialization of first label in ctor
#
# Appends 11 new recording rows (r=902..912) to the "DB" sheet for subject
# "024" / placement "BL-003_024_230528", and clears the saved "activeTab"
# override on the workbook so the first sheet is active by default again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

# New rows to append. Each entry is: recording timestamp (column D) and the
# protocol/label (column G). All other columns (A,B,C,E,F,H) are constant
# for this block of rows, matching the rows already above them.
$newRows = @(
    @("230528_131721", "Exp#1"),
    @("230528_132300", "Exp#1"),
    @("230528_132914", "Short#1_6s"),
    @("230528_133103", "Short#1_6s"),
    @("230528_133250", "Short#1_6s"),
    @("230528_133501", "Short#1_6s"),
    @("230528_133716", "Short#1_6s"),
    @("230528_133812", "Short#1_6s"),
    @("230528_134242", "Short#1_6s"),
    @("230528_134616", "Short#1_6s"),
    @("230528_135602", "Short#1_6s")
)

$firstRow = 902
$lastRow = $firstRow + $newRows.Count - 1

# Column A holds the zero-padded subject id "024". Typing that directly into
# a General-formatted cell would make Excel reinterpret it as the number 24,
# so it is written into a scratch range that is pre-formatted as Text, then
# copied as values into the real destination cells (this is the standard
# Excel trick for preserving leading zeros), after which the scratch range is
# cleared again.
$scratchA = $ws.Range("Z1:Z" + $newRows.Count)
$scratchA.NumberFormat = "@"
$scratchA.Value = "024"
$scratchA.Copy() | Out-Null
$ws.Range("A" + $firstRow + ":A" + $lastRow).PasteSpecial(-4163) | Out-Null
$scratchA.Clear() | Out-Null

# Fill in columns B:H one row at a time.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstRow + $i
    $timestamp = $newRows[$i][0]
    $label = $newRows[$i][1]

    $rowValues = New-Object 'object[,]' 1,7
    $rowValues[0,0] = "BL-003"
    $rowValues[0,1] = "BL-003_024_230528"
    $rowValues[0,2] = $timestamp
    $rowValues[0,3] = "1.1, 1.2"
    $rowValues[0,4] = "above knee"
    $rowValues[0,5] = $label
    $rowValues[0,6] = "left leg"

    $ws.Range("B" + $r + ":H" + $r).Value = $rowValues
}

# The workbook previously pinned the view to the 4th tab ("DB",
# activeTab="3"). Re-activate the first sheet so the saved view no longer
# overrides the default active tab.
$wb.Worksheets.Item(1).Activate()

Write-Host "done"
